$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Rename original sheet to reflect 31-count dataset
$ws1.Name = "4_pTHg_Comb_31ct"

# Duplicate the sheet (placed right after) to become the 30-count dataset
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "4_pTHg_Comb_30ct"

# The 30ct sheet excludes the anomalous low-flow reading (row 21, 2016-12-... 0.6 cfs)
$ws2.Rows(21).Delete()

# Flag the excluded row on the 31ct sheet with the highlight fill used for it
$ws1.Range("A21:D21").Interior.Color = 16249003

# Restore/adjust selections per the saved view state
$ws1.Range("D39").Select()
$ws2.Range("E32").Select()

# Make the 30ct sheet the active tab (matches activeTab="1" / tabSelected on sheet2)
$ws2.Activate()

Write-Host "done"
